$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")
$row = 70

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($row, 1).Value = 69

    # Date column: force Text format first so the ISO date string is not
    # auto-converted into a date serial number by Excel's input parser.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "12:54:07"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.06

    # Exit Price column: open trade has no exit price yet (empty string).
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.159115331514
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason column: open trade has no exit reason yet (empty string).
    $ws.Cells.Item($row, 16).NumberFormat = "@"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 16).Style = "Normal"

    $ws.Cells.Item($row, 17).Value = 0
}
